$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "price_coverage"
$ws.Columns.Item(1).ColumnWidth = 15.166666666666666
$ws.Range("B7").Select()
